$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 10.11875365514329
$ws.Range("D2").Value = 5.31023102338783
$ws.Range("E2").Value = 12.73687827976545
$ws.Range("F2").Value = 26.78117815595823
$ws.Range("G2").Value = 3.644883695392751
$ws.Range("I2").Value = 24.80523879970791
$ws.Range("K2").Value = 15.07934249687616
$ws.Range("L2").Value = 9.484170272418412
$ws.Range("N2").Value = 17.17579014067983
$ws.Range("O2").Value = 23.9150189516587

$ws.Range("C3").Value = 10.04089717265047
$ws.Range("D3").Value = 5.283706321572266
$ws.Range("E3").Value = 12.68178901720742
$ws.Range("F3").Value = 26.78629359557181
$ws.Range("G3").Value = 3.647100261956267
$ws.Range("I3").Value = 24.84200628154851
$ws.Range("K3").Value = 14.61568049680329
$ws.Range("L3").Value = 9.485601943859953
$ws.Range("N3").Value = 17.20846622824546
$ws.Range("O3").Value = 23.97460878926559

$ws.Range("C4").Value = 9.994826936652858
$ws.Range("D4").Value = 5.267177928899393
$ws.Range("E4").Value = 12.65061302324652
$ws.Range("F4").Value = 26.79754894776432
$ws.Range("G4").Value = 3.648533456826168
$ws.Range("I4").Value = 24.87095072505485
$ws.Range("K4").Value = 14.32488394323915
$ws.Range("L4").Value = 9.488113285934045
$ws.Range("N4").Value = 17.23035312883355
$ws.Range("O4").Value = 24.01740337559844

$ws.Range("C5").Value = 9.976506169989545
$ws.Range("D5").Value = 5.260383705310496
$ws.Range("E5").Value = 12.63858442200723
$ws.Range("F5").Value = 26.80417292021667
$ws.Range("G5").Value = 3.649135712655205
$ws.Range("I5").Value = 24.88434344238586
$ws.Range("K5").Value = 14.20502096920297
$ws.Range("L5").Value = 9.489547613399473
$ws.Range("N5").Value = 17.23973149931314
$ws.Range("O5").Value = 24.03639776165399

$ws.Range("C6").Value = 9.973491881773109
$ws.Range("D6").Value = 5.259252041815976
$ws.Range("E6").Value = 12.63662817034772
$ws.Range("F6").Value = 26.80539575525286
$ws.Range("G6").Value = 3.649236818776548
$ws.Range("I6").Value = 24.88666365826423
$ws.Range("K6").Value = 14.18504096836531
$ws.Range("L6").Value = 9.489810612884911
$ws.Range("N6").Value = 17.24131653093262
$ws.Range("O6").Value = 24.03964553397252

$ws.Range("C7").Value = 9.994577999129824
$ws.Range("D7").Value = 5.267086534800929
$ws.Range("E7").Value = 12.65044805252487
$ws.Range("F7").Value = 26.79763003717465
$ws.Range("G7").Value = 3.648541505220399
$ws.Range("I7").Value = 24.87112488090211
$ws.Range("K7").Value = 14.32327269976645
$ws.Range("L7").Value = 9.488130965393545
$ws.Range("N7").Value = 17.23047774827261
$ws.Range("O7").Value = 24.01765325079904

$ws.Range("C8").Value = 10.09155941820936
$ws.Range("D8").Value = 5.301136823105066
$ws.Range("E8").Value = 12.71733980483867
$ws.Range("F8").Value = 26.78125636587076
$ws.Range("G8").Value = 3.645633014779311
$ws.Range("I8").Value = 24.81659214697844
$ws.Range("K8").Value = 14.92083971732057
$ws.Range("L8").Value = 9.48432546435687
$ws.Range("N8").Value = 17.18667878069453
$ws.Range("O8").Value = 23.93427432984075

$ws.Range("C9").Value = 10.29470646383429
$ws.Range("D9").Value = 5.365896766478439
$ws.Range("E9").Value = 12.86906569465617
$ws.Range("F9").Value = 26.8136067543619
$ws.Range("G9").Value = 3.640499791039361
$ws.Range("I9").Value = 24.76033941020971
$ws.Range("K9").Value = 16.03723423969189
$ws.Range("L9").Value = 9.489787792361518
$ws.Range("N9").Value = 17.11522834022682
$ws.Range("O9").Value = 23.82023517271913

$ws.Range("C10").Value = 10.45076018473219
$ws.Range("D10").Value = 5.412128491315872
$ws.Range("E10").Value = 12.992387201418
$ws.Range("F10").Value = 26.87667727785051
$ws.Range("G10").Value = 3.637072355693501
$ws.Range("I10").Value = 24.75007689347242
$ws.Range("K10").Value = 16.81526167511391
$ws.Range("L10").Value = 9.501636659019942
$ws.Range("N10").Value = 17.07149572466827
$ws.Range("O10").Value = 23.76688802790579

$ws.Range("C11").Value = 10.52298512735475
$ws.Range("D11").Value = 5.43284246752793
$ws.Range("E11").Value = 13.0509016032057
$ws.Range("F11").Value = 26.91387230377234
$ws.Range("G11").Value = 3.635587015728075
$ws.Range("I11").Value = 24.75217530709401
$ws.Range("K11").Value = 17.15853782174018
$ws.Range("L11").Value = 9.508715794062061
$ws.Range("N11").Value = 17.05349479367205
$ws.Range("O11").Value = 23.74927915389329

$ws.Range("C12").Value = 10.5504914246686
$ws.Range("D12").Value = 5.440638533880198
$ws.Range("E12").Value = 13.07339216582969
$ws.Range("F12").Value = 26.9291745461653
$ws.Range("G12").Value = 3.635035110686875
$ws.Range("I12").Value = 24.75394343889097
$ws.Range("K12").Value = 17.28688709980657
$ws.Range("L12").Value = 9.511637969035556
$ws.Range("N12").Value = 17.04694988745905
$ws.Range("O12").Value = 23.74357203213544

$ws.Range("C13").Value = 10.54456084253331
$ws.Range("D13").Value = 5.438961676428703
$ws.Range("E13").Value = 13.06853387804627
$ws.Range("F13").Value = 26.92582490105552
$ws.Range("G13").Value = 3.635153504464191
$ws.Range("I13").Value = 24.75351934428784
$ws.Range("K13").Value = 17.25931944002298
$ws.Range("L13").Value = 9.510997911826099
$ws.Range("N13").Value = 17.04834737753774
$ws.Range("O13").Value = 23.74475837439904

$ws.Range("C14").Value = 10.52524507781443
$ws.Range("D14").Value = 5.433484826891081
$ws.Range("E14").Value = 13.05274533991831
$ws.Range("F14").Value = 26.91510686979168
$ws.Range("G14").Value = 3.635541398817794
$ws.Range("I14").Value = 24.75230126204898
$ws.Range("K14").Value = 17.16913068981535
$ws.Range("L14").Value = 9.508951372685676
$ws.Range("N14").Value = 17.05295089929148
$ws.Range("O14").Value = 23.74879034844175

$ws.Range("C15").Value = 10.51343332971149
$ws.Range("D15").Value = 5.430123797510759
$ws.Range("E15").Value = 13.04311724662576
$ws.Range("F15").Value = 26.9087001005627
$ws.Range("G15").Value = 3.63578036906668
$ws.Range("I15").Value = 24.75168193080853
$ws.Range("K15").Value = 17.11367058053125
$ws.Range("L15").Value = 9.507729208500438
$ws.Range("N15").Value = 17.05580604937689
$ws.Range("O15").Value = 23.75138528896452

$ws.Range("C16").Value = 10.44606315092481
$ws.Range("D16").Value = 5.410768209433373
$ws.Range("E16").Value = 12.98861048681717
$ws.Range("F16").Value = 26.87441717842324
$ws.Range("G16").Value = 3.637170906933852
$ws.Range("I16").Value = 24.75007597930939
$ws.Range("K16").Value = 16.79260372806779
$ws.Range("L16").Value = 9.501207891044174
$ws.Range("N16").Value = 17.07271017426338
$ws.Range("O16").Value = 23.76817310979865

$ws.Range("C17").Value = 10.40503498733761
$ws.Range("D17").Value = 5.398811353072433
$ws.Range("E17").Value = 12.95578049703243
$ws.Range("F17").Value = 26.85556034263677
$ws.Range("G17").Value = 3.638042824838265
$ws.Range("I17").Value = 24.75082449792284
$ws.Range("K17").Value = 16.59282906781609
$ws.Range("L17").Value = 9.497638919438838
$ws.Range("N17").Value = 17.08356479275761
$ws.Range("O17").Value = 23.78018002504643

$ws.Range("C18").Value = 10.38155424032192
$ws.Range("D18").Value = 5.391904376826863
$ws.Range("E18").Value = 12.9371257774622
$ws.Range("F18").Value = 26.8455152570622
$ws.Range("G18").Value = 3.638551280209999
$ws.Range("I18").Value = 24.75189197945947
$ws.Range("K18").Value = 16.47692710539542
$ws.Range("L18").Value = 9.495745216552281
$ws.Range("N18").Value = 17.0899863194295
$ws.Range("O18").Value = 23.78771269073944

$ws.Range("C19").Value = 10.37362490296188
$ws.Range("D19").Value = 5.38956076033731
$ws.Range("E19").Value = 12.93084924911807
$ws.Range("F19").Value = 26.8422518559795
$ws.Range("G19").Value = 3.63872462999671
$ws.Range("I19").Value = 24.7523627817853
$ws.Range("K19").Value = 16.4375170338104
$ws.Range("L19").Value = 9.495131403009186
$ws.Range("N19").Value = 17.09219116856258
$ws.Range("O19").Value = 23.79037062931946

$ws.Range("C20").Value = 10.40939049177541
$ws.Range("D20").Value = 5.400087272084305
$ws.Range("E20").Value = 12.95925179083613
$ws.Range("F20").Value = 26.85748483322589
$ws.Range("G20").Value = 3.637949288651257
$ws.Range("I20").Value = 24.75067888889302
$ws.Range("K20").Value = 16.61419942867644
$ws.Range("L20").Value = 9.498002388605483
$ws.Range("N20").Value = 17.0823908574987
$ws.Range("O20").Value = 23.77883699062148

$ws.Range("C21").Value = 10.53091451575192
$ws.Range("D21").Value = 5.435094827496456
$ws.Range("E21").Value = 13.05737391526194
$ws.Range("F21").Value = 26.9182220308614
$ws.Range("G21").Value = 3.635427178649495
$ws.Range("I21").Value = 24.7526326219396
$ws.Range("K21").Value = 17.1956666928084
$ws.Range("L21").Value = 9.509545949859264
$ws.Range("N21").Value = 17.05159136540162
$ws.Range("O21").Value = 23.74757995466562

$ws.Range("C22").Value = 10.61123841116743
$ws.Range("D22").Value = 5.457694240331073
$ws.Range("E22").Value = 13.12343277500735
$ws.Range("F22").Value = 26.96500888059417
$ws.Range("G22").Value = 3.633840366641209
$ws.Range("I22").Value = 24.7595835045716
$ws.Range("K22").Value = 17.56607577166822
$ws.Range("L22").Value = 9.518496967703481
$ws.Range("N22").Value = 17.03304526222697
$ws.Range("O22").Value = 23.73275419202553

$ws.Range("C23").Value = 10.5682926931052
$ws.Range("D23").Value = 5.445658879864327
$ws.Range("E23").Value = 13.0880043312979
$ws.Range("F23").Value = 26.93939123967863
$ws.Range("G23").Value = 3.634681665535254
$ws.Range("I23").Value = 24.7553545714683
$ws.Range("K23").Value = 17.36929445114758
$ws.Range("L23").Value = 9.513591448001012
$ws.Range("N23").Value = 17.04279900319379
$ws.Range("O23").Value = 23.74015335756093

$ws.Range("C24").Value = 10.40742103365806
$ws.Range("D24").Value = 5.399510531141404
$ws.Range("E24").Value = 12.9576817335271
$ws.Range("F24").Value = 26.85661229122526
$ws.Range("G24").Value = 3.637991553997429
$ws.Range("I24").Value = 24.75074273406693
$ws.Range("K24").Value = 16.60454114695043
$ws.Range("L24").Value = 9.497837571481055
$ws.Range("N24").Value = 17.08292102957205
$ws.Range("O24").Value = 23.77944221513438

$ws.Range("C25").Value = 10.23847340385909
$ws.Range("D25").Value = 5.348605391886747
$ws.Range("E25").Value = 12.82588677483757
$ws.Range("F25").Value = 26.79794733827075
$ws.Range("G25").Value = 3.641827795453516
$ws.Range("I25").Value = 24.77011099482248
$ws.Range("K25").Value = 16.03723423969189
$ws.Range("L25").Value = 9.486930226057193
$ws.Range("N25").Value = 17.13301597021071
$ws.Range("O25").Value = 23.84575921894702
